$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily draw result as row 13, mirroring the existing rows
# (all values stored as literal text, matching the "t=str" cells already
# present in rows 1-12). A leading apostrophe forces Excel to keep
# number/date-looking values ("2025-09-29", "250929",
# "2025-09-29T21:36:04.353+04:00") as plain text instead of auto-converting
# them to a date serial / number.
$ws.Range("A13").Value = "'2025-09-29"
$ws.Range("B13").Value = "Pick 4"
$ws.Range("C13").Value = "'250929"
$ws.Range("D13").Value = "2-7-2-0"
$ws.Range("E13").Value = "'2025-09-29T21:36:04.353+04:00"

# The sheet's dimension (A1:E12 -> A1:E13) is recalculated automatically by
# Excel once the new row is populated.

# Extend the "numbers stored as text" ignored-error suppression down to the
# new row, so row 13 behaves the same as rows 1-12 (no green
# number-as-text warning triangles).
$ws.Range("A1:E13").Errors.Item(9).Ignore = $true
